$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new row of homework data for 2020-02-13 at row 77
$row = 77

# Columns B, C, D hold text values in this sheet (date string, id string,
# ticker). Force text formatting before assigning so Excel doesn't
# auto-convert the date-like / numeric-like strings into a date serial or
# a number, then reset the style back to Normal so no extra style index is
# left on the cell (matches the unstyled data rows already in the sheet).
$ws.Range("B$row").NumberFormat = "@"
$ws.Range("C$row").NumberFormat = "@"
$ws.Range("D$row").NumberFormat = "@"

$ws.Cells.Item($row, 1).Value = 1581552000
$ws.Cells.Item($row, 2).Value = "2020-02-13"
$ws.Cells.Item($row, 3).Value = "5293"
$ws.Cells.Item($row, 4).Value = "AME"
$ws.Cells.Item($row, 5).Value = 1.75
$ws.Cells.Item($row, 6).Value = 1.93
$ws.Cells.Item($row, 7).Value = 1.72
$ws.Cells.Item($row, 8).Value = 1.93
$ws.Cells.Item($row, 9).Value = 4051400

$ws.Range("B$row").Style = "Normal"
$ws.Range("C$row").Style = "Normal"
$ws.Range("D$row").Style = "Normal"
